# Update thesis manuscript: add two note rows above the existing table and
# append a new block of observations (wavelet-based depth/residency checks)
# below the existing data, matching the author's latest notes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two blank rows at the very top; this pushes the existing
#     header + data (old rows 1-14) down to rows 3-16, carrying their
#     styles (bold header row) and values along with them. ---
$ws.Rows("1:2").Insert()

# --- Fill the two new note rows (column A only) ---
$ws.Range("A1").Value = "have  7-10 day rolling window "
$ws.Range("A2").Value = "maybe also have a 30-50 days window to look at longer periods, for the residency? Eg depth median change/day wavelet has a significant bereich over ~100 days"

# --- Append the new observation rows (17-28) below the existing table ---
$newRows = @(
    @("wm and sm", "m", 321, "median depth wavelet", "high period frequencies in periods between 2 - 10 days"),
    @("sr", "m", 321, "median depth wavelet", "low period frequencies in periods between 2 - 15 days"),
    @("wm", "m", 321, "median depth change wavelet (roll3)", "high period frequencies in periods either around 2 days or 32 days"),
    @("sm", "m", 321, "median depth change wavelet (roll3)", "high period frequencies in periods between 2 and 15 days AND low frequencies around 64 and 128 days"),
    @("wm", "m", 321, "depth range wavelet", "significant high frequencies > 15 days"),
    @("wm", "f", 308, "median depth wavelet", "significant high frequencies between 3 and 8 days"),
    @("wr", "f", 308, "median depth wavelet", "significant high frequencies between 8 and 16 days"),
    @("sr", "f", 308, "median depth wavelet", "low period frequencies in periods between 2 - 15 days (dark blue)"),
    @("wr", "f", 308, "median depth change wavelet (roll3)", "significant high frequencies between 2 and 17 days"),
    @("sr", "f", 308, "median depth change wavelet (roll3)", "low period frequencies in periods between 2 - 15 days (yellow)"),
    @("wm", "f", 308, "median depth change wavelet (roll3)", "significant high frequencies between 3 and 8 days"),
    @("sr", "f", 308, "depth range wavelet", "low period frequencies in periods between 2 - 15 days (blue-green)")
)

$r = 17
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

# --- Column width tweaks: column A widened for the longer note text,
#     column D widened for the longer wavelet-category labels. The COM
#     ColumnWidth setter here only resolves to the nearest 1/6-character
#     step, so we pick the inputs that land closest to the authored
#     widths (28.09 / 31.27 characters). ---
$ws.Columns("A").ColumnWidth = 27.3
$ws.Columns("D").ColumnWidth = 30.5

# --- Update the saved selection / scroll position to match where the
#     author was last working. ---
$ws.Range("E17").Select()
